# 自动更新Excel文件 - 2026-01-25 23:13:53
# For every data row (row 2..99), decrement the "剩余" (remaining days, column E)
# by 1. When remaining would drop to 0, the cycle restarts: remaining resets back
# to the row's "总天" (total days, column D) and the "开始时间" (start date,
# column F) is stamped with the new cycle's start date (20260126).
#
# Row 36 holds a corrupted start-date value (202510929) and is left untouched,
# matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCycleStart = 20260126
$lastRow = 99

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) {
        continue
    }

    $totalDays = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2

    if ($remaining -eq $null -or $totalDays -eq $null) {
        continue
    }

    $newRemaining = $remaining - 1

    if ($newRemaining -eq 0) {
        $newRemaining = $totalDays
        $ws.Cells.Item($r, 6).Value = $newCycleStart
    }

    $ws.Cells.Item($r, 5).Value = $newRemaining
}
